$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 20.447252
$ws.Range("H2").Value = 61.341756
$ws.Range("I2").Value = 0.8699145605694745
$ws.Range("J2").Value = 0.8770588936480435
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 16.60495066666667
$ws.Range("N2").Value = 49.814852
$ws.Range("O2").Value = 0.4330603147186406
$ws.Range("P2").Value = 0.5197056776409935
$ws.Range("Q2").Value = 339.5256107289013
$ws.Range("R2").Value = 3055.730496560112
$ws.Range("S2").Value = 0.3767254733785446
$ws.Range("T2").Value = 0.4558124866544165
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 20.447252
$ws.Range("H3").Value = 61.341756
$ws.Range("I3").Value = 0.8699145605694745
$ws.Range("J3").Value = 0.8770588936480435
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 0.7472513333333333
$ws.Range("N3").Value = 2.241754
$ws.Range("O3").Value = 0.01948845883877707
$ws.Range("P3").Value = 0.02338764916283215
$ws.Range("Q3").Value = 15.27923632000267
$ws.Range("R3").Value = 137.513126880024
$ws.Range("S3").Value = 0.01695329410691104
$ws.Range("T3").Value = 0.02051234569978216
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 20.447252
$ws.Range("H4").Value = 61.341756
$ws.Range("I4").Value = 0.8699145605694745
$ws.Range("J4").Value = 0.8770588936480435
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 0.5008306666666666
$ws.Range("N4").Value = 1.502492
$ws.Range("O4").Value = 0.01306176034372721
$ws.Range("P4").Value = 0.01567511679067463
$ws.Range("Q4").Value = 10.24061085066133
$ws.Range("R4").Value = 92.165497655952
$ws.Range("S4").Value = 0.01136261550967724
$ws.Range("T4").Value = 0.01374800059023296
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 20.447252
$ws.Range("H5").Value = 61.341756
$ws.Range("I5").Value = 0.8699145605694745
$ws.Range("J5").Value = 0.8770588936480435
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 1.312462666666667
$ws.Range("N5").Value = 3.937388
$ws.Range("O5").Value = 0.034229279381366
$ws.Range("P5").Value = 0.04107776730272161
$ws.Range("Q5").Value = 26.83625488592533
$ws.Range("R5").Value = 241.526293973328
$ws.Range("S5").Value = 0.02977654853165078
$ws.Range("T5").Value = 0.03602762114405678
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 20.447252
$ws.Range("H6").Value = 61.341756
$ws.Range("I6").Value = 0.8699145605694745
$ws.Range("J6").Value = 0.8770588936480435
$ws.Range("K6").Value = 2
$ws.Range("M6").Value = 19.1777795
$ws.Range("N6").Value = 38.355559
$ws.Range("O6").Value = 0.5001601867174891
$ws.Range("P6").Value = 0.4001537891027781
$ws.Range("Q6").Value = 392.132890236934
$ws.Range("R6").Value = 2352.797341421604
$ws.Range("S6").Value = 0.4350966290426908
$ws.Range("T6").Value = 0.3509584395595551
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 2.483247333333333
$ws.Range("H7").Value = 7.449742
$ws.Range("I7").Value = 0.105648084777455
$ws.Range("J7").Value = 0.1065157390747562
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 16.60495066666667
$ws.Range("N7").Value = 49.814852
$ws.Range("O7").Value = 0.4330603147186406
$ws.Range("P7").Value = 0.5197056776409935
$ws.Range("Q7").Value = 41.23419946313156
$ws.Range("R7").Value = 371.107795168184
$ws.Range("S7").Value = 0.04575199284314629
$ws.Range("T7").Value = 0.05535683435527745
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 2.483247333333333
$ws.Range("H8").Value = 7.449742
$ws.Range("I8").Value = 0.105648084777455
$ws.Range("J8").Value = 0.1065157390747562
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 0.7472513333333333
$ws.Range("N8").Value = 2.241754
$ws.Range("O8").Value = 0.01948845883877707
$ws.Range("P8").Value = 0.02338764916283215
$ws.Range("Q8").Value = 1.855609880829778
$ws.Range("R8").Value = 16.700488927468
$ws.Range("S8").Value = 0.002058918351581062
$ws.Range("T8").Value = 0.002491152735800171
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 2.483247333333333
$ws.Range("H9").Value = 7.449742
$ws.Range("I9").Value = 0.105648084777455
$ws.Range("J9").Value = 0.1065157390747562
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 0.5008306666666666
$ws.Range("N9").Value = 1.502492
$ws.Range("O9").Value = 0.01306176034372721
$ws.Range("P9").Value = 0.01567511679067463
$ws.Range("Q9").Value = 1.243686417451556
$ws.Range("R9").Value = 11.193177757064
$ws.Range("S9").Value = 0.001379949964136892
$ws.Range("T9").Value = 0.001669646650041829
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 2.483247333333333
$ws.Range("H10").Value = 7.449742
$ws.Range("I10").Value = 0.105648084777455
$ws.Range("J10").Value = 0.1065157390747562
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.312462666666667
$ws.Range("N10").Value = 3.937388
$ws.Range("O10").Value = 0.034229279381366
$ws.Range("P10").Value = 0.04107776730272161
$ws.Range("Q10").Value = 3.259169417099555
$ws.Range("R10").Value = 29.332524753896
$ws.Range("S10").Value = 0.003616257809953748
$ws.Range("T10").Value = 0.004375428743790248
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 2.483247333333333
$ws.Range("H11").Value = 7.449742
$ws.Range("I11").Value = 0.105648084777455
$ws.Range("J11").Value = 0.1065157390747562
$ws.Range("K11").Value = 2
$ws.Range("M11").Value = 19.1777795
$ws.Range("N11").Value = 38.355559
$ws.Range("O11").Value = 0.5001601867174891
$ws.Range("P11").Value = 0.4001537891027781
$ws.Range("Q11").Value = 47.62316980262967
$ws.Range("R11").Value = 285.739018815778
$ws.Range("S11").Value = 0.052840965808637
$ws.Range("T11").Value = 0.04262267658984655
$ws.Range("E12").Value = 2
$ws.Range("G12").Value = 0.5743975
$ws.Range("H12").Value = 1.148795
$ws.Range("I12").Value = 0.02443735465307048
$ws.Range("J12").Value = 0.01642536727720028
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 16.60495066666667
$ws.Range("N12").Value = 49.814852
$ws.Range("O12").Value = 0.4330603147186406
$ws.Range("P12").Value = 0.5197056776409935
$ws.Range("Q12").Value = 9.537842150556667
$ws.Range("R12").Value = 57.22705290334
$ws.Range("S12").Value = 0.01058284849694974
$ws.Range("T12").Value = 0.008536356631299574
$ws.Range("E13").Value = 2
$ws.Range("G13").Value = 0.5743975
$ws.Range("H13").Value = 1.148795
$ws.Range("I13").Value = 0.02443735465307048
$ws.Range("J13").Value = 0.01642536727720028
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 0.7472513333333333
$ws.Range("N13").Value = 2.241754
$ws.Range("O13").Value = 0.01948845883877707
$ws.Range("P13").Value = 0.02338764916283215
$ws.Range("Q13").Value = 0.4292192977383333
$ws.Range("R13").Value = 2.57531578643
$ws.Range("S13").Value = 0.0004762463802849614
$ws.Range("T13").Value = 0.0003841507272498238
$ws.Range("E14").Value = 2
$ws.Range("G14").Value = 0.5743975
$ws.Range("H14").Value = 1.148795
$ws.Range("I14").Value = 0.02443735465307048
$ws.Range("J14").Value = 0.01642536727720028
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 0.5008306666666666
$ws.Range("N14").Value = 1.502492
$ws.Range("O14").Value = 0.01306176034372721
$ws.Range("P14").Value = 0.01567511679067463
$ws.Range("Q14").Value = 0.2876758828566667
$ws.Range("R14").Value = 1.72605529714
$ws.Range("S14").Value = 0.0003191948699130736
$ws.Range("T14").Value = 0.0002574695503998397
$ws.Range("E15").Value = 2
$ws.Range("G15").Value = 0.5743975
$ws.Range("H15").Value = 1.148795
$ws.Range("I15").Value = 0.02443735465307048
$ws.Range("J15").Value = 0.01642536727720028
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.312462666666667
$ws.Range("N15").Value = 3.937388
$ws.Range("O15").Value = 0.034229279381366
$ws.Range("P15").Value = 0.04107776730272161
$ws.Range("Q15").Value = 0.7538752745766666
$ws.Range("R15").Value = 4.52325164746
$ws.Range("S15").Value = 0.0008364730397614741
$ws.Range("T15").Value = 0.0008364730397614741
$ws.Range("E16").Value = 2
$ws.Range("G16").Value = 0.5743975
$ws.Range("H16").Value = 1.148795
$ws.Range("I16").Value = 0.02443735465307048
$ws.Range("J16").Value = 0.01642536727720028
$ws.Range("K16").Value = 2
$ws.Range("M16").Value = 19.1777795
$ws.Range("N16").Value = 38.355559
$ws.Range("O16").Value = 0.5001601867174891
$ws.Range("P16").Value = 0.4001537891027781
$ws.Range("Q16").Value = 11.01566860035125
$ws.Range("R16").Value = 44.062674401405
$ws.Range("S16").Value = 0.01222259186616123
$ws.Range("T16").Value = 0.006572672953376475
